$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 309, shifting existing rows 309-378 down to 310-379.
$ws.Rows.Item(309).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(309, 1).Value = 5
$ws.Cells.Item(309, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(309, 3).Value = "Maule"
$ws.Cells.Item(309, 4).Value = 45204
$ws.Cells.Item(309, 5).Value = 7
$ws.Cells.Item(309, 6).Value = 100112017
$ws.Cells.Item(309, 7).Value = "Apio"
$ws.Cells.Item(309, 8).Value = "Americana (o)"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 700
$ws.Cells.Item(309, 11).Value = 6000
$ws.Cells.Item(309, 12).Value = 6000
$ws.Cells.Item(309, 13).Value = 6000
$ws.Cells.Item(309, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(309, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(309, 16).Value = 500
$ws.Cells.Item(309, 17).Value = 12
$ws.Cells.Item(309, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Cells.Item(309, 4).NumberFormat = $ws.Cells.Item(310, 4).NumberFormat
